# Update performance document: add a new "v1289" run's column (G) of
# pipeline-performance numbers to the "Sponza" and "ComplexMesh" sheets.
# The "PartOfSponza" sheet is left untouched.

$wb = $excel.ActiveWorkbook

$sponza = $wb.Worksheets.Item("Sponza")
$complexMesh = $wb.Worksheets.Item("ComplexMesh")

# --- New column header (shared string "v1289") --------------------------
$sponza.Range("G1").Value = "v1289"
$complexMesh.Range("G1").Value = "v1289"

# --- New per-run values for "Sponza" -------------------------------------
$sponzaValues = @(10176, 10206, 10194, 10139, 10157, 10184, 10201, 10175, 10125, 10187)
for ($i = 0; $i -lt $sponzaValues.Length; $i++) {
    $row = 2 + $i
    $sponza.Cells.Item($row, 7).Value = $sponzaValues[$i]
}

# --- New per-run values for "ComplexMesh" --------------------------------
$complexMeshValues = @(7683, 7657, 7648, 7612, 7716, 7690, 7631, 7615, 7645, 7643)
for ($i = 0; $i -lt $complexMeshValues.Length; $i++) {
    $row = 2 + $i
    $complexMesh.Cells.Item($row, 7).Value = $complexMeshValues[$i]
}

# The AVG/VAR/DIFF-ACCEPT/ratio rows (12-16) already contain formulas that
# span through column G (shared formulas / per-cell formulas referencing
# G2:G11, G12 etc.) - they recompute automatically now that G has data,
# turning the previous #DIV/0! results into real numbers.

# --- Update the selections left on each sheet ----------------------------
# Visit "Sponza" first so its own selection is preserved, then finish on
# "ComplexMesh" so it remains the active/selected tab (as in the source).
$sponza.Range("G15").Select() | Out-Null
$complexMesh.Range("G15").Select() | Out-Null
